$wb = $excel.ActiveWorkbook

# Update status text "Ready for handoff" -> "In Translation" on all sheets that reference it
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "In Translation"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "In Translation"

# Narrow the status columns (Overview E:F, zh-cn C, de-de C) from 17.2159881591797 to 13.4101848602295
$wsOverview.Range("E:F").ColumnWidth = 12.5

$wsZh.Range("C:C").ColumnWidth = 12.5
$wsDe.Range("C:C").ColumnWidth = 12.5
